$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transformed facilities list")

# Rows 60-102, column B: "Olin E. Teague Veterans Center" -> insert " before Center
for ($r = 60; $r -le 102; $r++) {
    $ws.Cells.Item($r, 2).Value = "Olin E. Teague Veterans`" Center"
}

# Rows 298-337, column B: "Audie L. Murphy Memorial Veterans Hospital" -> insert " before Hospital
for ($r = 298; $r -le 337; $r++) {
    $ws.Cells.Item($r, 2).Value = "Audie L. Murphy Memorial Veterans`" Hospital"
}

# Rows 542-570, column B: "George H. OBrien, Jr. ..." -> insert " between O and Brien
for ($r = 542; $r -le 570; $r++) {
    $ws.Cells.Item($r, 2).Value = "George H. O`"Brien, Jr. Department of Veterans Affairs Medical Center"
}
